# Weekly fruit/vegetable price update:
# Insert a new row at 113 (pushing the previous rows 113-114 down to 114-115)
# and populate it with the new week's "Caigua" / "Primera" price record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing row 113 (and everything below it) down by one row.
$ws.Rows.Item(113).Insert()

# Fill in the newly inserted row 113 with this week's data.
$ws.Cells.Item(113, 1).Value = 1
$ws.Cells.Item(113, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(113, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(113, 4).Value = 44753
$ws.Cells.Item(113, 5).Value = 15
$ws.Cells.Item(113, 6).Value = 100112036
$ws.Cells.Item(113, 7).Value = "Caigua"
$ws.Cells.Item(113, 8).Value = "Sin especificar"
$ws.Cells.Item(113, 9).Value = "Primera"
$ws.Cells.Item(113, 10).Value = 120
$ws.Cells.Item(113, 11).Value = 6000
$ws.Cells.Item(113, 12).Value = 7000
$ws.Cells.Item(113, 13).Value = 6500
$ws.Cells.Item(113, 14).Value = "$/caja 20 kilos"
$ws.Cells.Item(113, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(113, 16).Value = 325
$ws.Cells.Item(113, 17).Value = 20
$ws.Cells.Item(113, 18).Value = "Hortaliza"
